# Insert a new data row before the current row 19 ("Vega Monumental
# Concepción" / Espinaca weekly price sheet), pushing all rows from 19
# downward (old 19..96 -> new 20..97) and populate the newly inserted
# row with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 19:96 down by one row to make room for the new record.
$ws.Rows(19).Insert()

# Populate the new row 19 with the new weekly record.
$ws.Range("A19").Value = 11
$ws.Range("B19").Value = "Vega Monumental Concepción"
$ws.Range("C19").Value = "Bíobío"
$ws.Range("D19").Value = 44881
$ws.Range("E19").Value = 8
$ws.Range("F19").Value = 100112012
$ws.Range("G19").Value = "Espinaca"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 50
$ws.Range("K19").Value = 13000
$ws.Range("L19").Value = 14000
$ws.Range("M19").Value = 13400
$ws.Range("N19").Value = "$/cuna 10 kilos"
$ws.Range("O19").Value = "Región Metropolitana"
$ws.Range("P19").Value = 1340
$ws.Range("Q19").Value = 10
$ws.Range("R19").Value = "Hortaliza"
